$wb = $excel.ActiveWorkbook

# --- "Worksheet" sheet: German (column C) and French (column F) translations ---
$ws1 = $wb.Worksheets.Item("Worksheet")

$ws1.Range("C2").Value = "Willkommen!"
$ws1.Range("F2").Value = "вітаємо вас щиро!"

$ws1.Range("C3").Value = "Speichern als..."
$ws1.Range("F3").Value = "Enregistrer sous..."

$ws1.Range("C4").Value = "Ansicht"
$ws1.Range("F4").Value = "cava"

$ws1.Range("C5").Value = "Über..."
$ws1.Range("F5").Value = "À propos..."

$ws1.Range("C6").Value = "Über mich"
$ws1.Range("F6").Value = "À propos de moi"

$ws1.Range("C7").Value = "Aktuelles Passwort"
$ws1.Range("F7").Value = "Mot de passe actuel"

$ws1.Range("C8").Value = "Neues Passwort"

$ws1.Range("C9").Value = "Neues Passwort bestätigen"

$ws1.Range("C10").Value = "Passwort ändern"

$ws1.Range("C13").Value = "Nachricht lesen"

$ws1.Range("C14").Value = "Sind Sie sicher, dass Sie diese Nachricht löschen möchten?"

$ws1.Range("C15").Value = "In Nachrichten suchen"

$ws1.Range("C16").Value = "Nachricht verfassen"

$ws1.Range("C17").Value = "Kein Datum angegeben"

$ws1.Range("C18").Value = "Schnellstart"

# --- "Sheet1" sheet: French (column C) translations ---
$ws2 = $wb.Worksheets.Item("Sheet1")

$ws2.Range("C2").Value = " Traduction 1"
$ws2.Range("F2").Value = " Traducción 1"

$ws2.Range("C3").Value = "Traduction 2"

$ws2.Range("C4").Value = "Traduction 3 "

$ws2.Range("C5").Value = "Traduction 4"

$ws2.Range("C6").Value = "Traduction 5"
